# Weekly fruit/vegetable price update: a new week's record is inserted at
# the top of the data block (row 48), pushing the existing rows 48:90 down
# to 49:91 (the sheet's "latest first" log gains one more entry).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 48; existing rows 48-90 shift down to 49-91.
$ws.Rows("48:48").Insert()

# Populate the newly inserted row 48 with this week's record.
$ws.Range("A48").Value = 3
$ws.Range("B48").Value = "Femacal de La Calera"
$ws.Range("C48").Value = "Coquimbo"
$ws.Range("D48").Value = 44827
$ws.Range("E48").Value = 5
$ws.Range("F48").Value = 100112035
$ws.Range("G48").Value = "Bruselas (repollito)"
$ws.Range("H48").Value = "Sin especificar"
$ws.Range("I48").Value = "Primera"
$ws.Range("J48").Value = 45
$ws.Range("K48").Value = 15000
$ws.Range("L48").Value = 15000
$ws.Range("M48").Value = 15000
$ws.Range("N48").Value = "$/malla 15 kilos"
$ws.Range("O48").Value = "Provincia de Quillota"
$ws.Range("P48").Value = 1000
$ws.Range("Q48").Value = 15
$ws.Range("R48").Value = "Hortaliza"
